$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("training")

$row = 24

$ws.Cells.Item($row, 1).Value = "2025-09-16 13:29:08"
$ws.Cells.Item($row, 2).Value = "training"
$ws.Cells.Item($row, 3).Value = "configs/training/2025-09-16/z/0000"

$noteCell = $ws.Cells.Item($row, 4)
$noteCell.Value = "'"
$noteCell.ClearFormats()

$ws.Cells.Item($row, 5).Value = "['cross_entropy']"
$ws.Cells.Item($row, 6).Value = "[1.0]"
$ws.Cells.Item($row, 7).Value = "['torch.optim.adam.Adam']"
$ws.Cells.Item($row, 8).Value = "[0.001]"
$ws.Cells.Item($row, 9).Value = 128
$ws.Cells.Item($row, 10).Value = 128
$ws.Cells.Item($row, 11).Value = "general_utils.ml.training.NoImprovementStopping"
$ws.Cells.Item($row, 12).Value = 8
$ws.Cells.Item($row, 13).Value = 0.00001
$ws.Cells.Item($row, 14).Value = 500
